$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A12").Value = "en reporte de OT resaltar observaciones"
$ws.Range("B12").Value = "no comenzado"

$ws.Range("C10").Select()
